$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row (row 1) updates ---
# C1: "max" -> "prediction"
$ws.Range("C1").Value = "prediction"
# D1 used to hold "prediction"; it now takes what used to be in E1: "rejection-f"
$ws.Range("D1").Value = "rejection-f"

# --- Row 2 updates (this row keeps the former row-3 data: even_MAG-GUT57690.fa) ---
$ws.Range("A2").Value = "even_MAG-GUT57690.fa"
$ws.Range("B2").Value = -15.76543825665451
$ws.Range("C2").Value = "s__CAG-791 sp000431495"
$ws.Range("D2").Value = "s__CAG-791 sp000431495(reject)"

# --- Remove now-unused column E entirely (was rejection-f / duplicate data) ---
$ws.Columns.Item(5).Delete()

# --- Remove rows 3 through 6 (only one data row remains) ---
$ws.Range("A3:A6").EntireRow.Delete()
